# Edit: rename "总计" sheet (sheetId=6) to "2022-Q1" and populate it with the
# new quarterly fund-holdings breakdown; add a fresh sheet (becomes sheetId=7)
# named "总计" right after it, carrying forward the summary table with the new
# 2022-Q1 row inserted at the top.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: the existing 6th sheet (currently "总计") becomes "2022-Q1".
# ---------------------------------------------------------------------------
$q1Sheet = $wb.Worksheets.Item(6)
$q1Sheet.Name = "2022-Q1"

# Wipe its old "总计" summary content (A1:D6) before laying out the new table.
$q1Sheet.Cells.Clear()

# ---------------------------------------------------------------------------
# Step 2: add a brand-new sheet right after "2022-Q1" — this will carry the
# "总计" name forward and gets a fresh sheetId.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1Sheet)
$totalSheet.Name = "总计"

# ---------------------------------------------------------------------------
# Step 3: header row + data for the "2022-Q1" fund-holdings sheet.
# ---------------------------------------------------------------------------
$q1Headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $q1Headers.Length; $c++) {
    $q1Sheet.Cells.Item(1, 2 + $c).Value = $q1Headers[$c]
}

$fundData = @(
    @(0, "002910", "易方达供给改革灵活配置混合", "49.29", "87.54", "8.25", "4.0664", 2),
    @(1, "001373", "易方达新丝路灵活配置混合", "38.95", "90.79", "3.65", "1.4217", 7),
    @(2, "001225", "中邮趋势精选灵活配置混合", "12.03", "70.17", "5.36", "0.6448", 6),
    @(3, "000925", "汇添富外延增长主题股票A", "20.62", "83.57", "3.04", "0.6268", 5),
    @(4, "012008", "易方达稳健回报一年封闭运作混合型证券投资基金A", "18.94", "50.67", "2.47", "0.4678", 9),
    @(5, "959991", "兴证资管金麒麟领先优势一年持有期混合A", "8.28", "79.70", "4.79", "0.3966", 3),
    @(6, "012175", "易方达稳健增利混合型证券投资基金A", "13.94", "53.92", "2.58", "0.3597", 10),
    @(7, "011777", "易方达稳健增长混合型证券投资基金A", "13.88", "54.12", "2.58", "0.3581", 8),
    @(8, "000545", "中邮核心竞争力灵活配置混合", "3.43", "73.83", "7.24", "0.2483", 2),
    @(9, "011346", "淳厚鑫淳一年持有期混合型证券投资基金", "5.75", "67.80", "3.69", "0.2122", 2),
    @(10, "012454", "淳厚鑫悦混合A", "3.31", "76.84", "3.87", "0.1281", 2),
    @(11, "009488", "中邮价值精选混合A", "1.32", "78.12", "8.11", "0.1071", 2),
    @(12, "001484", "天弘新价值灵活配置混合", "3.51", "89.24", "2.74", "0.0962", 6),
    @(13, "012009", "易方达稳健回报一年封闭运作混合型证券投资基金C", "2.03", "50.67", "2.47", "0.0501", 9),
    @(14, "582003", "东吴配置优化灵活配置混合", "1.04", "90.74", "4.26", "0.0443", 4),
    @(15, "002281", "建信裕利灵活配置混合", "1.10", "88.94", "3.64", "0.0400", 6),
    @(16, "002378", "建信弘利灵活配置混合", "1.03", "89.57", "3.66", "0.0377", 7),
    @(17, "008846", "大成民稳增长混合A", "2.49", "22.89", "1.50", "0.0374", 3),
    @(18, "009489", "中邮价值精选混合C", "0.39", "78.12", "8.11", "0.0316", 2),
    @(19, "012455", "淳厚鑫悦混合C", "0.79", "76.84", "3.87", "0.0306", 2),
    @(20, "007254", "广发均衡价值混合", "0.49", "89.66", "5.45", "0.0267", 5),
    @(21, "959993", "兴证资管金麒麟领先优势一年持有期混合C", "0.43", "79.70", "4.79", "0.0206", 3),
    @(22, "012176", "易方达稳健增利混合型证券投资基金C", "0.66", "53.92", "2.58", "0.0170", 10),
    @(23, "006225", "人保量化基本面混合A", "0.63", "88.00", "1.74", "0.0110", 3),
    @(24, "011778", "易方达稳健增长混合型证券投资基金C", "0.38", "54.12", "2.58", "0.0098", 8),
    @(25, "008847", "大成民稳增长混合C", "0.56", "22.89", "1.50", "0.0084", 3),
    @(26, "009796", "大成汇享一年持有期混合A", "0.38", "22.99", "1.46", "0.0055", 4),
    @(27, "011424", "汇添富外延增长主题股票C", "0.07", "83.57", "3.04", "0.0021", 5),
    @(28, "006226", "人保量化基本面混合C", "0.04", "88.00", "1.74", "0.0007", 3),
    @(29, "009797", "大成汇享一年持有期混合C", "0.04", "22.99", "1.46", "0.0006", 4)
)


# Force the numeric-looking text columns (fund code / scale / position /
# weight / value) to be stored as text, matching the source data (these are
# strings in the original table, not numbers).
$q1Sheet.Range("B2:B31").NumberFormat = "@"
$q1Sheet.Range("D2:G31").NumberFormat = "@"

for ($i = 0; $i -lt $fundData.Length; $i++) {
    $row = $fundData[$i]
    $r = 2 + $i
    $q1Sheet.Cells.Item($r, 1).Value = $row[0]
    $q1Sheet.Cells.Item($r, 2).Value = $row[1]
    $q1Sheet.Cells.Item($r, 3).Value = $row[2]
    $q1Sheet.Cells.Item($r, 4).Value = $row[3]
    $q1Sheet.Cells.Item($r, 5).Value = $row[4]
    $q1Sheet.Cells.Item($r, 6).Value = $row[5]
    $q1Sheet.Cells.Item($r, 7).Value = $row[6]
    $q1Sheet.Cells.Item($r, 8).Value = $row[7]
}

# Match the house style used on every other quarter sheet: bold/bordered
# header row and bold/bordered index column (style index 2 in styles.xml).
$sourceStyle = $wb.Worksheets.Item(5)
$sourceStyle.Range("B1:H1").Copy()
$q1Sheet.Range("B1:H1").PasteSpecial(-4122)
$sourceStyle.Range("A2:A31").Copy()
$q1Sheet.Range("A2:A31").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Step 4: header row + data for the "总计" (grand total) sheet.
# ---------------------------------------------------------------------------
$totalHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($c = 0; $c -lt $totalHeaders.Length; $c++) {
    $totalSheet.Cells.Item(1, 2 + $c).Value = $totalHeaders[$c]
}

$totalData = @(
    @(0, "2022-Q1", 30, 9.51),
    @(1, "2021-Q4", 41, 8.970000000000001),
    @(2, "2021-Q3", 10, 2.95),
    @(3, "2021-Q2", 10, 1.82),
    @(4, "2021-Q1", 3, 1.13),
    @(5, "2020-Q4", 1, 0.05)
)


for ($i = 0; $i -lt $totalData.Length; $i++) {
    $row = $totalData[$i]
    $r = 2 + $i
    $totalSheet.Cells.Item($r, 1).Value = $row[0]
    $totalSheet.Cells.Item($r, 2).Value = $row[1]
    $totalSheet.Cells.Item($r, 3).Value = $row[2]
    $totalSheet.Cells.Item($r, 4).Value = $row[3]
}

$sourceStyle.Range("B1:D1").Copy()
$totalSheet.Range("B1:D1").PasteSpecial(-4122)
$sourceStyle.Range("A2:A7").Copy()
$totalSheet.Range("A2:A7").PasteSpecial(-4122)
